# Applies the recomputed Profit-tracker figures for Sheets/Malboro_Profits.xlsx
# (commit: "chore: update Sheets via scheduled runner").
# Each block below targets one worksheet; cell addresses/values were
# derived from the canonical OOXML diff of the commit.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1337193.2
$ws.Range("I17").Value = 1034.6154
$ws.Range("K17").Value = 3103.8462
$ws.Range("M17").Value = -2935.8462
# Row 113
$ws.Range("H113").Value = 14005.826
$ws.Range("I113").Value = 14478.177
$ws.Range("J113").Value = 12667.5
$ws.Range("K113").Value = 14478.177
$ws.Range("L113").Value = 12667.5
$ws.Range("M113").Value = -11224.177
$ws.Range("N113").Value = -19175.5
# Row 132
$ws.Range("H132").Value = 8102.2827
$ws.Range("I132").Value = 6311.5
$ws.Range("K132").Value = 18934.5
$ws.Range("M132").Value = -16404.5
# Row 133
$ws.Range("H133").Value = 65500.332
$ws.Range("J133").Value = 65500.332
$ws.Range("L133").Value = 65500.332
$ws.Range("N133").Value = -75620.33199999999

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 11053.179
$ws.Range("I32").Value = 4127.1187
$ws.Range("K32").Value = 4127.1187
$ws.Range("M32").Value = -3840.1187
# Row 74
$ws.Range("H74").Value = 17101.965
$ws.Range("I74").Value = 2455.6667
$ws.Range("K74").Value = 2455.6667
$ws.Range("M74").Value = -1581.6667
# Row 77
$ws.Range("H77").Value = 17101.965
$ws.Range("I77").Value = 2455.6667
$ws.Range("K77").Value = 12278.3335
$ws.Range("M77").Value = -7910.333500000001
# Row 122
$ws.Range("H122").Value = 2901.8333
$ws.Range("I122").Value = 1283.1666
$ws.Range("K122").Value = 3849.4998
$ws.Range("M122").Value = -1399.4998
# Row 132
$ws.Range("H132").Value = 2950218.5
$ws.Range("I132").Value = 1260
$ws.Range("K132").Value = 3780
$ws.Range("M132").Value = -1250

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 24999
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 24999
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 24999
$ws.Range("N26").Value = -25583
$ws.Range("M26").ClearContents()
# Row 134
$ws.Range("H134").Value = 14177.655
$ws.Range("J134").Value = 21536.23
$ws.Range("L134").Value = 64608.69
$ws.Range("N134").Value = -69678.69

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 13976.361
$ws.Range("I58").Value = 7332.4287
$ws.Range("J58").Value = 18204.318
$ws.Range("K58").Value = 7332.4287
$ws.Range("L58").Value = 18204.318
$ws.Range("M58").Value = -7129.4287
$ws.Range("N58").Value = -18610.318
# Row 86
$ws.Range("H86").Value = 11012.643
$ws.Range("J86").Value = 8998
$ws.Range("L86").Value = 8998
$ws.Range("N86").Value = -11244
# Row 89
$ws.Range("H89").Value = 11012.643
$ws.Range("J89").Value = 8998
$ws.Range("L89").Value = 44990
$ws.Range("N89").Value = -56222
# Row 105
$ws.Range("H105").Value = 17367.555
$ws.Range("I105").Value = 25577.25
$ws.Range("K105").Value = 25577.25
$ws.Range("M105").Value = -23830.25
# Row 112
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42954
# Row 136
$ws.Range("H136").Value = 13976.361
$ws.Range("I136").Value = 7332.4287
$ws.Range("J136").Value = 18204.318
$ws.Range("K136").Value = 21997.2861
$ws.Range("L136").Value = 54612.954
$ws.Range("M136").Value = -19447.2861
$ws.Range("N136").Value = -59712.954

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 47.5
$ws.Range("I12").Value = 70.5
$ws.Range("K12").Value = 211.5
$ws.Range("M12").Value = -38.5
# Row 92
$ws.Range("H92").Value = 914.5
$ws.Range("J92").Value = 818.125
$ws.Range("L92").Value = 2454.375
$ws.Range("N92").Value = -4950.375
# Row 116
$ws.Range("H116").Value = 3618.7778
$ws.Range("I116").Value = 3943
$ws.Range("K116").Value = 11829
$ws.Range("M116").Value = -8387
# Row 131
$ws.Range("H131").Value = 1484.09
$ws.Range("J131").Value = 1494.0303
$ws.Range("L131").Value = 4482.090899999999
$ws.Range("N131").Value = -14562.0909

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 13097.667
$ws.Range("I80").Value = 5625.3335
$ws.Range("J80").Value = 23060.777
$ws.Range("K80").Value = 5625.3335
$ws.Range("L80").Value = 23060.777
$ws.Range("M80").Value = -4627.3335
$ws.Range("N80").Value = -25056.777
# Row 83
$ws.Range("H83").Value = 13097.667
$ws.Range("I83").Value = 5625.3335
$ws.Range("J83").Value = 23060.777
$ws.Range("K83").Value = 28126.6675
$ws.Range("L83").Value = 115303.885
$ws.Range("M83").Value = -23134.6675
$ws.Range("N83").Value = -125287.885
# Row 102
$ws.Range("H102").Value = 1651.2858
$ws.Range("I102").Value = 1481.7916
$ws.Range("K102").Value = 1481.7916
$ws.Range("M102").Value = 140.2084
# Row 122
$ws.Range("H122").Value = 1996.25
$ws.Range("I122").Value = 1748.75
$ws.Range("K122").Value = 5246.25
$ws.Range("M122").Value = -2796.25
# Row 132
$ws.Range("H132").Value = 17524.389
$ws.Range("I132").Value = 12914.0625
$ws.Range("K132").Value = 38742.1875
$ws.Range("M132").Value = -36212.1875

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 3252915.8
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 3252915.8
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3252915.8
$ws.Range("N68").Value = -3254413.8
$ws.Range("M68").ClearContents()
# Row 71
$ws.Range("H71").Value = 3252915.8
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 3252915.8
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 16264579
$ws.Range("N71").Value = -16272067
$ws.Range("M71").ClearContents()
# Row 82
$ws.Range("H82").Value = 5158
$ws.Range("I82").Value = 4389.6
$ws.Range("K82").Value = 4389.6
$ws.Range("M82").Value = -4028.6
# Row 85
$ws.Range("H85").Value = 5158
$ws.Range("I85").Value = 4389.6
$ws.Range("K85").Value = 4389.6
$ws.Range("M85").Value = -3141.6
# Row 122
$ws.Range("H122").Value = 7825.645
$ws.Range("I122").Value = 5606.067
$ws.Range("J122").Value = 9906.5
$ws.Range("K122").Value = 16818.201
$ws.Range("L122").Value = 29719.5
$ws.Range("M122").Value = -14368.201
$ws.Range("N122").Value = -34619.5
# Row 132
$ws.Range("H132").Value = 2235516
$ws.Range("I132").Value = 8424.25
$ws.Range("K132").Value = 25272.75
$ws.Range("M132").Value = -22742.75

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 69
$ws.Range("H69").Value = 58360.5
$ws.Range("J69").Value = 58360.5
$ws.Range("L69").Value = 58360.5
$ws.Range("N69").Value = -59858.5
# Row 72
$ws.Range("H72").Value = 58360.5
$ws.Range("J72").Value = 58360.5
$ws.Range("L72").Value = 175081.5
$ws.Range("N72").Value = -182569.5
# Row 100
$ws.Range("H100").Value = 1422.3
$ws.Range("I100").Value = 1509.25
$ws.Range("K100").Value = 3018.5
$ws.Range("M100").Value = -2477.5
# Row 122
$ws.Range("H122").Value = 3785.1836
$ws.Range("I122").Value = 1918.5
$ws.Range("K122").Value = 5755.5
$ws.Range("M122").Value = -3305.5

